# Rename the inline picture shapes living in the document's headers and
# footers:
#   - the two "Pearson logo" pictures in the footers (both currently named
#     "image1.png") become "image2.png"
#   - the "BTec logo" picture in the first-page header (currently named
#     "image2.jpg") becomes "image1.jpg"

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {

    # --- Footers: Pearson logo pictures: image1.png -> image2.png
    foreach ($ftr in $sec.Footers) {
        if ($ftr.Exists) {
            foreach ($ishp in $ftr.Range.InlineShapes) {
                if ($ishp.AlternativeText -like "*PearsonLogo.png") {
                    $ishp.Select() | Out-Null
                    $word.Selection.InlineShapes.Item(1).Name = "image2.png"
                }
            }
        }
    }

    # --- Headers: BTec logo picture: image2.jpg -> image1.jpg
    foreach ($hdr in $sec.Headers) {
        if ($hdr.Exists) {
            foreach ($ishp in $hdr.Range.InlineShapes) {
                if ($ishp.AlternativeText -eq "BTec_Logo-Orange") {
                    $ishp.Select() | Out-Null
                    $word.Selection.InlineShapes.Item(1).Name = "image1.jpg"
                }
            }
        }
    }
}
